$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the end of the Prog5.c
#    paragraph. In the edited document it has moved to sit inside the
#    third "make clean&&..." command line further down, so drop it here;
#    it will be re-created (embedded in the inserted XML) below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Append the large block of new paragraphs (Prog11-14, Prog21-24,
#    the Command: header and the three "make clean&&..." lines) after
#    the existing trailing empty paragraph. We use Range.InsertXML with
#    a full WordprocessingML package fragment so that we can faithfully
#    reproduce the exact run/paragraph layout (including the eastAsia
#    rFonts hints, the custom indentation on the two "send a message"
#    continuation lines, and the gramStart/gramEnd proofing marks)
#    rather than relying on typed-text heuristics.
# ------------------------------------------------------------------
$newParasXml = @'
    <w:p/>
    <w:p>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>rog</w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>.c: listen port 5, send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>2.c (port 6)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>2.c: listen port 6, send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">.c (port </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>3.c:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>listen port 7,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>14</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">.c (port </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>8</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>4.c:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>listen port 8,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>13</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">.c (port </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>7</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>rog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>.c: listen port 5, send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>2.c (port 6)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="2160"/>
      </w:pPr>
      <w:r>
        <w:t>send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">.c (port </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>7</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="2160"/>
      </w:pPr>
      <w:r>
        <w:t>send a message to prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">.c (port </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>8</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t>2.c: listen port 6</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t>3.c:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>listen port 7</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Prog</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t>4.c:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>listen port 8</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>C</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>ommand</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>make clean&amp;&amp;make depend&amp;&amp;make&amp;</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>&amp;./</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>nachos -x ../test/prog1 -x ../test/prog2 -x ../test/prog3 -x ../test/prog4 -x ../test/prog5</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>make clean&amp;&amp;make depend&amp;&amp;make&amp;</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>&amp;./</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>nachos -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>1 -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>2 -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>3 -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">4 </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>make clean&amp;&amp;make depend&amp;&amp;make&amp;</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>&amp;./</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>nachos -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>2</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>1 -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t>2 -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t>3 -x ../test/prog</w:t>
      </w:r>
      <w:r>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">4 </w:t>
      </w:r>
    </w:p>
'@

$pkgXml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $newParasXml + '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

# The very last (currently empty) paragraph in the document is the
# InsertXML merge target: its identity (paraId etc.) is kept, but its
# content becomes the last paragraph of $newParasXml ("...prog24 "). The
# two blank trailing paragraphs that follow it in the final document are
# then appended as brand-new empty paragraphs.
$insertAt = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$insertAt.InsertXML($pkgXml)

$tail = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$tail.InsertParagraphAfter()
$tail2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$tail2.InsertParagraphAfter()

Write-Output "done"
